# Update the two-digit x two-digit multiplication problems in the table.
# The source table has 25 problem cells laid out as 5 "content" rows
# (rows 1, 5, 10, 15, 20 of the 20-row table) x 5 columns each, with
# blank spacer rows in between. Some source problems repeat verbatim
# (e.g. "66x46=" appears twice), so cells are addressed positionally
# via Table.Cell(row, col) rather than via a document-wide Find/Replace,
# which would otherwise replace every matching occurrence at once.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    @(1, 1, "51×84="),
    @(1, 2, "51×45="),
    @(1, 3, "53×90="),
    @(1, 4, "57×11="),
    @(1, 5, "81×84="),

    @(5, 1, "26×97="),
    @(5, 2, "96×75="),
    @(5, 3, "24×57="),
    @(5, 4, "84×48="),
    @(5, 5, "21×35="),

    @(10, 1, "54×15="),
    @(10, 2, "39×49="),
    @(10, 3, "20×65="),
    @(10, 4, "91×97="),
    @(10, 5, "30×37="),

    @(15, 1, "89×80="),
    @(15, 2, "50×12="),
    @(15, 3, "54×96="),
    @(15, 4, "35×51="),
    @(15, 5, "15×87="),

    @(20, 1, "17×45="),
    @(20, 2, "53×47="),
    @(20, 3, "45×97="),
    @(20, 4, "27×27="),
    @(20, 5, "87×42=")
)

foreach ($entry in $newValues) {
    $row = $entry[0]
    $col = $entry[1]
    $text = $entry[2]
    $t.Cell($row, $col).Range.Text = $text
}
